# Actualización automática del mapa (mapa_interactivo_AYKO.html)
#
# A new incident report row (Caso -500, Castañares 5656) was inserted
# at the top of this block, pushing the existing rows down by one and
# dropping the last (oldest) row off the bottom of the table.
#
# Net effect on the sheet:
#   row 89 -> new data (Caso -500)
#   row 90 -> what used to be row 89 (Caso 6394)
#   row 91 -> what used to be row 90 (Caso 6399)
#   (what used to be row 91, Caso -515, falls off the table)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that must stay as TEXT even though their contents look numeric
# (Caso, F. De Reclamo, Comuna, OT) - force text format before writing so
# Excel doesn't silently convert them to numbers / dates.
$textCols = @("A", "B", "C", "D", "E", "F", "G", "H", "J", "K", "L", "O", "P")

function Set-RowData {
    param(
        $row,
        $caso,
        $fReclamo,
        $direccion,
        $comuna,
        $ot,
        $proveedor,
        $estado,
        $observaciones,
        $attachments,
        $tipoTarea,
        $equipo,
        $tipoElemento,
        $coordX,
        $coordY,
        $operacion,
        $zona
    )

    foreach ($col in $textCols) {
        $ws.Range("$col$row").NumberFormat = "@"
    }

    $ws.Range("A$row").Value = $caso
    $ws.Range("B$row").Value = $fReclamo
    $ws.Range("C$row").Value = $direccion
    $ws.Range("D$row").Value = $comuna
    $ws.Range("E$row").Value = $ot
    $ws.Range("F$row").Value = $proveedor
    $ws.Range("G$row").Value = $estado
    $ws.Range("H$row").Value = $observaciones
    $ws.Range("I$row").Value = $attachments
    $ws.Range("J$row").Value = $tipoTarea
    $ws.Range("K$row").Value = $equipo
    $ws.Range("L$row").Value = $tipoElemento
    $ws.Range("M$row").Value = $coordX
    $ws.Range("N$row").Value = $coordY
    $ws.Range("O$row").Value = $operacion
    $ws.Range("P$row").Value = $zona
}

# Row 89: brand-new case that was inserted
Set-RowData 89 "-500" "7/3/2025" "Castañares 5656" "8" "807965768" `
    "AYKO" "Pendiente" "Columna chocada con rienda a pique" 1 `
    "Cambio" "Sin equipos" "Terminal" -58.479921 -34.673021 `
    "Boedo" "Capital Sur"

# Row 90: former row 89 data, shifted down
Set-RowData 90 "6394" "7/14/2025" "LAMBARE 1076" "5" "808194286" `
    "AYKO" "Pendiente" "Picada" 1 `
    "Cambio" "Sin equipos" "Pasante" -58.43008 -34.601416 `
    "Almagro" "Capital Sur"

# Row 91: former row 90 data, shifted down (former row 91 data is dropped)
Set-RowData 91 "6399" "7/14/2025" "ESCALADA AV. 966" "9" "808258198" `
    "AYKO" "Pendiente" "Picada" 1 `
    "Cambio" "Sin equipos" "Pasante" -58.493069 -34.646557 `
    "Devoto" "Capital Norte"

$wb.Save()
